$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New note added for the final week (row 26, column D) describing the
#    last week's work: building the management-screen layout and the post
#    management feature. This introduces a new shared string and a new
#    highlighted fill for the cell.
# ---------------------------------------------------------------------------
$ws.Range("D26").Value = "Tuần cuối dựng layout cho màn hình quản lý, làm chức năng quản lý các bài post"
# "Green, Accent 6, Lighter 40%" background highlight for the new note.
$ws.Range("D26").Interior.Color = 9359785

# ---------------------------------------------------------------------------
# 2) Center-align the "Thời gian chi tiết" (detailed date) column, header and
#    all the date cells below it, to make the table look tidier / more
#    responsive.
# ---------------------------------------------------------------------------
$ws.Range("B1").HorizontalAlignment = -4108

$dateRows = @(2,3,4,5,6,7,8,10,11,12,13,14,15,17,18,19,20,21,22,24,25,26)
foreach ($r in $dateRows) {
    $ws.Cells.Item($r, 2).HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 3) Resize the columns for a more responsive-looking layout: narrower date
#    columns, much wider content columns so long notes need fewer wrapped
#    lines.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 11.88
$ws.Columns.Item(2).ColumnWidth = 14.88
$ws.Columns.Item(3).ColumnWidth = 64.59
$ws.Columns.Item(4).ColumnWidth = 83.74

# Row 5 wraps less now that column D is much wider.
$ws.Rows.Item(5).RowHeight = 60

# ---------------------------------------------------------------------------
# 4) Update the view so that the sheet opens scrolled near the bottom of the
#    plan with the last note selected.
# ---------------------------------------------------------------------------
$ws.Range("C24").Select()

Write-Host "done"
